$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.148.43"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.558.09"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.38"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.779.71"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.552.16"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.86"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.155.02"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.32"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.95"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.56"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0461"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.379.39"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.945"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.810"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.512"
$ws.Range("E40").Value = "  -4.55%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.22"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.692.60"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.40"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  +0.03%  "
